# Append four new "Ref" field rows (63-66) to the format sheet, mirroring
# the existing "field" rows for exchanges / parameters / production volume /
# properties, each documenting a new mandatory "Ref" field of type "str".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
    @{ Row = 63; A = "exchanges";          B = "Ref"; C = "str"; E = "mandatory"; G = "Ref" },
    @{ Row = 64; A = "parameters";         B = "Ref"; C = "str"; E = "mandatory"; G = "Ref" },
    @{ Row = 65; A = "production volume";  B = "Ref"; C = "str"; E = "mandatory"; G = "Ref" },
    @{ Row = 66; A = "properties";         B = "Ref"; C = "str"; E = "mandatory"; G = "Ref" }
)

foreach ($r in $newRows) {
    $ws.Cells.Item($r.Row, 1).Value = $r.A   # A: parent
    $ws.Cells.Item($r.Row, 2).Value = $r.B   # B: field name
    $ws.Cells.Item($r.Row, 3).Value = $r.C   # C: format
    $ws.Cells.Item($r.Row, 5).Value = $r.E   # E: mandatory/optional
    $ws.Cells.Item($r.Row, 7).Value = $r.G   # G: in dataframe
}

# Match the new selection / scroll state from the commit: the active cell
# block is the newly added column-A rows.
$ws.Range("A63:A66").Select() | Out-Null
